$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old first header row ("...mation", "...pompes)", Hiver/Eté/Année) —
# this shifts the old units header (row 2) and the data rows (old 3-7) up by one,
# so the data now lives in rows 2-6 and the sheet shrinks to A1:K6.
$ws.Rows("1").Delete()

# Rebuild row 1 as a single header row with the new column labels
# (idx / idx2 / Name / Date Start / Date End / units...).
$ws.Range("A1").Value = "idx"
$ws.Range("B1").Value = "idx2"
$ws.Range("C1").Value = "Name"
$ws.Range("D1").Value = "Date Start"
$ws.Range("E1").Value = "Date End"
$ws.Range("F1").Value = "(m3/s)"
$ws.Range("G1").Value = "(MW1)"
$ws.Range("H1").Value = "(MW2)"
$ws.Range("I1").Value = "(GWh) Winter"
$ws.Range("J1").Value = "(GWh) Summer"
$ws.Range("K1").Value = "(GWh) Year"

# A1:E1 keep the default (unstyled) look; F1:K1 pick up the Arial 9 / General
# header style (a new cell style, distinct from the plain default one).
$ws.Range("F1:K1").Font.Name = "Arial"
$ws.Range("F1:K1").Font.Size = 9
$ws.Range("F1:K1").NumberFormatLocal = "General"

# The active selection moves to the first data row.
$ws.Range("A2:K2").Select()
